$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "Save" header in H1, matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the value for the new "Save" column in H2
$ws.Range("H2").Value = 0
